{"js": "// Helper: wrap a sequence of <w:p>...</w:p> (or run/bookmark fragments) into a\n// full OOXML package string acceptable to Range.insertOoxml().\nfunction wrapOoxml(bodyFragment) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyFragment + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\n// Common run-property blocks used throughout this section of the document.\nconst RPR_B = '<w:rPr><w:rFonts w:cstheme=\"minorHAnsi\"/><w:b/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr>';\nconst RPR = '<w:rPr><w:rFonts w:cstheme=\"minorHAnsi\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr>';\nconst PPR_B = '<w:pPr><w:rPr><w:rFonts w:cstheme=\"minorHAnsi\"/><w:b/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr></w:pPr>';\nconst PPR = '<w:pPr><w:rPr><w:rFonts w:cstheme=\"minorHAnsi\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr></w:pPr>';\n\nfunction t(text) {\n  const esc = text.replace(/&/g, '&amp;').replace(/</g, '&lt;').replace(/>/g, '&gt;');\n  const needsPreserve = /^\\s|\\s$/.test(text);\n  return needsPreserve ? '<w:t xml:space=\"preserve\">' + esc + '</w:t>' : '<w:t>' + esc + '</w:t>';\n}\n\nasync function replaceParagraphWhole(paragraph, innerXml) {\n  const range = paragraph.getRange('Whole');\n  range.insertOoxml(wrapOoxml('<w:p>' + innerXml + '</w:p>'), Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Locate all paragraphs once.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load('items/text');\nawait context.sync();\n\nfunction findParagraphIndex(exactText, startAt) {\n  for (let i = startAt || 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === exactText) return i;\n  }\n  throw new Error('paragraph not found: ' + exactText);\n}\n\nconst idxQ3 = findParagraphIndex('Q3. Docker vs Vagrant');\n\n// 1) \"Q3. Docker vs Vagrant\" -> \"Q3. \" + \"Docker \" + \"Container \" + bookmark(_GoBack) + \"Metrics\"\nawait replaceParagraphWhole(\n  paragraphs.items[idxQ3],\n  PPR_B +\n    '<w:r>' + RPR_B + '<w:lastRenderedPageBreak/>' + t('Q3. ') + '</w:r>' +\n    '<w:r>' + RPR_B + t('Docker ') + '</w:r>' +\n    '<w:r>' + RPR_B + t('Container ') + '</w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r>' + RPR_B + t('Metrics') + '</w:r>'\n);\n\n// 2) \"I used Apache JMeter ... docket. My findings are as follows:\" -> 5 runs\nconst idxApache = findParagraphIndex(\n  'I used Apache JMeter benchmarking tool to perform load testing on the virtual machine and used an extension called Container Watch for docket. My findings are as follows:'\n);\nawait replaceParagraphWhole(\n  paragraphs.items[idxApache],\n  PPR +\n    '<w:r>' + RPR + t('I used Apache JMeter benchmarking tool to perform load testing on ') + '</w:r>' +\n    '<w:r>' + RPR + t('HTTP requests ') + '</w:r>' +\n    '<w:r>' + RPR + t('and used an extension called Container Watch for ') + '</w:r>' +\n    '<w:r>' + RPR + t('CPU and Memory usage') + '</w:r>' +\n    '<w:r>' + RPR + t('. My findings are as follows:') + '</w:r>'\n);\n\n// 3) Delete the standalone \"Docker:\" paragraph entirely.\nconst idxDocker = findParagraphIndex('Docker:');\nparagraphs.items[idxDocker].delete();\nawait context.sync();\n\n// Re-fetch paragraphs since indices shift after the delete.\nparagraphs.load('items/text');\nawait context.sync();\n\n// 4) \"As the graph shows, CPU utilization ... 18% - 25%.\" -> 5 runs\nconst idxGraph = findParagraphIndex(\n  'As the graph shows, CPU utilization bumps up-to more than 60% when the request simulation is carried out. After the simulation is complete, it comes down in the range of 18% - 25%.'\n);\nawait replaceParagraphWhole(\n  paragraphs.items[idxGraph],\n  PPR +\n    '<w:r>' + RPR + t('As the graph shows, CPU ') + '</w:r>' +\n    '<w:r>' + RPR + t('and Memory ') + '</w:r>' +\n    '<w:r>' + RPR + t('utilization bumps up-to more than 60% when the request simulation is carried out. After the simulation is complete, ') + '</w:r>' +\n    '<w:r>' + RPR + t('memory utilization') + '</w:r>' +\n    '<w:r>' + RPR + t(' comes down in the range of 18% - 25%.') + '</w:r>'\n);\n\n// 5) \"JMeter:\" -> \"JMeter\" + \" Load Testing\" + \":\"\nparagraphs.load('items/text');\nawait context.sync();\nconst idxJMeter = findParagraphIndex('JMeter:');\nawait replaceParagraphWhole(\n  paragraphs.items[idxJMeter],\n  PPR_B +\n    '<w:r>' + RPR_B + t('JMeter') + '</w:r>' +\n    '<w:r>' + RPR_B + t(' Load Testing') + '</w:r>' +\n    '<w:r>' + RPR_B + t(':') + '</w:r>'\n);\n\n// 6) \"This table shows...\" paragraph: drop the _GoBack bookmark and the\n//    lastRenderedPageBreak marker (both move to the \"Error %:\" paragraph below).\nparagraphs.load('items/text');\nawait context.sync();\nconst idxTable = findParagraphIndex(\n  'This table shows the summary of the above 1000 HTTP requests. Some of the important metrics are:'\n);\nawait replaceParagraphWhole(\n  paragraphs.items[idxTable],\n  PPR_B +\n    '<w:r>' + RPR + t('This table shows the summary of the above 1000 HTTP requests. Some of the important metrics are:') + '</w:r>'\n);\n\n// 7) \"Error %: 0\" paragraph: add lastRenderedPageBreak to the first (\"Error %: \") run.\nparagraphs.load('items/text');\nawait context.sync();\nconst idxError = findParagraphIndex('Error %: 0');\nawait replaceParagraphWhole(\n  paragraphs.items[idxError],\n  PPR +\n    '<w:r>' + RPR_B + '<w:lastRenderedPageBreak/>' + t('Error %: ') + '</w:r>' +\n    '<w:r>' + RPR + '<w:t>0</w:t>' + '</w:r>'\n);\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# Helpers\n# ---------------------------------------------------------------------------\n\nfunction Find-ParagraphIndex($doc, [string]$exactText) {\n    $paras = $doc.Paragraphs\n    for ($i = 1; $i -le $paras.Count; $i++) {\n        $p = $paras.Item($i)\n        $txt = $p.Range.Text.TrimEnd([char]13, [char]7)\n        if ($txt -eq $exactText) {\n            return $i\n        }\n    }\n    throw (\"paragraph not found: \" + $exactText)\n}\n\nfunction Wrap-Ooxml([string]$bodyFragment) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $bodyFragment + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n$RPR_B = '<w:rPr><w:rFonts w:cstheme=\"minorHAnsi\"/><w:b/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr>'\n$RPR   = '<w:rPr><w:rFonts w:cstheme=\"minorHAnsi\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr>'\n$PPR_B = '<w:pPr><w:rPr><w:rFonts w:cstheme=\"minorHAnsi\"/><w:b/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr></w:pPr>'\n$PPR   = '<w:pPr><w:rPr><w:rFonts w:cstheme=\"minorHAnsi\"/><w:sz w:val=\"20\"/><w:szCs w:val=\"20\"/></w:rPr></w:pPr>'\n\nfunction Text-Run([string]$txt) {\n    $esc = $txt.Replace(\"&\", \"&amp;\").Replace(\"<\", \"&lt;\").Replace(\">\", \"&gt;\")\n    if ($txt -match '^\\s' -or $txt -match '\\s$') {\n        return '<w:t xml:space=\"preserve\">' + $esc + '</w:t>'\n    } else {\n        return '<w:t>' + $esc + '</w:t>'\n    }\n}\n\nfunction Replace-ParagraphWhole($doc, [string]$exactText, [string]$innerXml) {\n    $idx = Find-ParagraphIndex $doc $exactText\n    $p = $doc.Paragraphs.Item($idx)\n    $rng = $p.Range\n    $rng.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1; exclude the paragraph mark\n    $rng.Text = \"\"\n    $xml = Wrap-Ooxml ('<w:p>' + $innerXml + '</w:p>')\n    $rng.InsertXML($xml)\n}\n\n# ---------------------------------------------------------------------------\n# 1) \"Q3. Docker vs Vagrant\" -> \"Q3. \" + \"Docker \" + \"Container \" + bookmark(_GoBack) + \"Metrics\"\n# ---------------------------------------------------------------------------\n$inner1 = $PPR_B +\n    '<w:r>' + $RPR_B + '<w:lastRenderedPageBreak/>' + (Text-Run \"Q3. \") + '</w:r>' +\n    '<w:r>' + $RPR_B + (Text-Run \"Docker \") + '</w:r>' +\n    '<w:r>' + $RPR_B + (Text-Run \"Container \") + '</w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:r>' + $RPR_B + (Text-Run \"Metrics\") + '</w:r>'\nReplace-ParagraphWhole $d \"Q3. Docker vs Vagrant\" $inner1\n\n# ---------------------------------------------------------------------------\n# 2) \"I used Apache JMeter ... docket. My findings are as follows:\" -> 5 runs\n# ---------------------------------------------------------------------------\n$origApache = \"I used Apache JMeter benchmarking tool to perform load testing on the virtual machine and used an extension called Container Watch for docket. My findings are as follows:\"\n$inner2 = $PPR +\n    '<w:r>' + $RPR + (Text-Run \"I used Apache JMeter benchmarking tool to perform load testing on \") + '</w:r>' +\n    '<w:r>' + $RPR + (Text-Run \"HTTP requests \") + '</w:r>' +\n    '<w:r>' + $RPR + (Text-Run \"and used an extension called Container Watch for \") + '</w:r>' +\n    '<w:r>' + $RPR + (Text-Run \"CPU and Memory usage\") + '</w:r>' +\n    '<w:r>' + $RPR + (Text-Run \". My findings are as follows:\") + '</w:r>'\nReplace-ParagraphWhole $d $origApache $inner2\n\n# ---------------------------------------------------------------------------\n# 3) Delete the standalone \"Docker:\" paragraph entirely.\n# ---------------------------------------------------------------------------\n$idxDocker = Find-ParagraphIndex $d \"Docker:\"\n$d.Paragraphs.Item($idxDocker).Range.Delete()\n\n# ---------------------------------------------------------------------------\n# 4) \"As the graph shows, CPU utilization ... 18% - 25%.\" -> 5 runs\n# ---------------------------------------------------------------------------\n$origGraph = \"As the graph shows, CPU utilization bumps up-to more than 60% when the request simulation is carried out. After the simulation is complete, it comes down in the range of 18% - 25%.\"\n$inner4 = $PPR +\n    '<w:r>' + $RPR + (Text-Run \"As the graph shows, CPU \") + '</w:r>' +\n    '<w:r>' + $RPR + (Text-Run \"and Memory \") + '</w:r>' +\n    '<w:r>' + $RPR + (Text-Run \"utilization bumps up-to more than 60% when the request simulation is carried out. After the simulation is complete, \") + '</w:r>' +\n    '<w:r>' + $RPR + (Text-Run \"memory utilization\") + '</w:r>' +\n    '<w:r>' + $RPR + (Text-Run \" comes down in the range of 18% - 25%.\") + '</w:r>'\nReplace-ParagraphWhole $d $origGraph $inner4\n\n# ---------------------------------------------------------------------------\n# 5) \"JMeter:\" -> \"JMeter\" + \" Load Testing\" + \":\"\n# ---------------------------------------------------------------------------\n$inner5 = $PPR_B +\n    '<w:r>' + $RPR_B + (Text-Run \"JMeter\") + '</w:r>' +\n    '<w:r>' + $RPR_B + (Text-Run \" Load Testing\") + '</w:r>' +\n    '<w:r>' + $RPR_B + (Text-Run \":\") + '</w:r>'\nReplace-ParagraphWhole $d \"JMeter:\" $inner5\n\n# ---------------------------------------------------------------------------\n# 6) \"This table shows...\" paragraph: drop the _GoBack bookmark and the\n#    lastRenderedPageBreak marker (both move to the \"Error %:\" paragraph below).\n# ---------------------------------------------------------------------------\n$origTable = \"This table shows the summary of the above 1000 HTTP requests. Some of the important metrics are:\"\n$inner6 = $PPR_B +\n    '<w:r>' + $RPR + (Text-Run $origTable) + '</w:r>'\nReplace-ParagraphWhole $d $origTable $inner6\n\n# ---------------------------------------------------------------------------\n# 7) \"Error %: 0\" paragraph: add lastRenderedPageBreak to the first (\"Error %: \") run.\n# ---------------------------------------------------------------------------\n$inner7 = $PPR +\n    '<w:r>' + $RPR_B + '<w:lastRenderedPageBreak/>' + (Text-Run \"Error %: \") + '</w:r>' +\n    '<w:r>' + $RPR + '<w:t>0</w:t></w:r>'\nReplace-ParagraphWhole $d \"Error %: 0\" $inner7\n"}
